$d = $word.ActiveDocument

function Replace-Exact($oldText, $newText) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "Could not find text: $oldText"
    }
    $rng = $d.Range($idx, $idx + $oldText.Length)
    $rng.Text = $newText
}

# ---------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------
Replace-Exact "Digital Innovation Shaping the Future" "Chemistry and Its Profound Influence on the World Around Us"

# ---------------------------------------------------------------------
# 2. Author name: "Ashley Ward" -> "Dr. Isabella Sinclair"
# ---------------------------------------------------------------------
Replace-Exact "Ashley Ward" "Dr. Isabella Sinclair"

# ---------------------------------------------------------------------
# 3. Username / email line
# ---------------------------------------------------------------------
Replace-Exact "ward" "isabella"
Replace-Exact "ashley01@evansville" "sinclair@schooledulink"

# ---------------------------------------------------------------------
# 4. Body paragraph (sz 24) - four segments separated by <w:br/><w:br/>
# ---------------------------------------------------------------------

# Segment 1 (gains a fourth sentence)
$seg1old = "In an era of constant technological advancements, the impact of digital innovation on various aspects of modern civilization has been profound. Digital technologies, encompassing the internet, artificial intelligence, data analytics, robotics, and blockchain, have revolutionized communication, information sharing, and the overall functioning of society. Their transformative potential is evident in diverse fields, from the way we conduct business and engage in entertainment to the delivery of healthcare and governance."
$seg1new = "Step into the mesmerizing realm of Chemistry, where matter transforms and elements dance in harmony. At the heart of Chemistry lies the study of the composition, properties, and behavior of substances. Chemistry permeates every aspect of our lives, shaping the world we see, feel, taste, and smell. Dive into this extraordinary realm of science to unlock the secrets that unravel the intricate connections between matter and our existence."
Replace-Exact $seg1old $seg1new

# Segment 2
$seg2old = "In the realm of business, digital innovation has enabled companies to reach a global audience, streamline operations, and enhance customer experiences. E-commerce giants like Amazon and Alibaba have paved the way for seamless online shopping, while social media platforms like Facebook and Twitter have fostered unprecedented levels of connectivity and engagement. Digital tools have also spurred a surge in remote work, unlocking new possibilities for flexible employment and work-life balance."
$seg2new = "In our quest to understand the universe, Chemistry emerges as a vital link, cementing our comprehension of materials, energy, and reactions. This intricate field encompasses diverse areas spanning from the synthesis of novel pharmaceuticals that conquer formidable diseases to the creation of sustainable technologies that safeguard our planet. Chemistry wields the transformative power to shape our lives and propel us towards a brighter future."
Replace-Exact $seg2old $seg2new

# Segment 3
$seg3old = "In the healthcare sector, digital innovation is improving patient care and transforming the delivery of medical services. Telemedicine platforms provide remote consultations, allowing patients to access quality healthcare from the comfort of their homes. Wearable health devices monitor vital signs and share data with healthcare providers, enabling proactive health management. Big data analytics play a crucial role in processing vast amounts of medical data to identify patterns, predict outbreaks, and develop personalized treatment plans."
$seg3new = "Embark on a voyage through the annals of Chemistry, witnessing the remarkable achievements that have reshaped society. Discover how the discovery of elements like radium revolutionized medical treatments, leading to life-saving therapies. Experience the transformative nature of Chemistry in action, as synthetic materials revolutionized industries, reshaping communication, transportation, and the way we live. Marvel at the advancements in biotechnology, where genetic modifications enhance crop yields, promising to alleviate global hunger."
Replace-Exact $seg3old $seg3new

# Segment 4 (Culture and entertainment ...) is removed entirely, along with
# the two <w:br/> runs that introduced it.
$full = $d.Content.Text
$cultureIdx = $full.IndexOf("Culture and entertainment")
$delStart = $cultureIdx - 2   # back up over the two break characters
$delEnd = $full.IndexOf("massive revenues") + "massive revenues".Length + 1  # include trailing period
$delRange = $d.Range($delStart, $delEnd)
$delRange.Text = ""

# ---------------------------------------------------------------------
# 5. Summary body paragraph
# ---------------------------------------------------------------------
$sumOld = "Digital innovation has propelled the world into a new era, where technology permeates every aspect of life, transforming our ways of working, communicating, accessing entertainment, and receiving healthcare. These advancements offer immense potential for improving lives, yet they also raise questions about privacy, security, and the impact on employment. As the digital landscape continues to evolve, it is imperative to navigate its complexities and harness its transformative power for the benefit of society."
$sumNew = "Chemistry unveils the complexity of matter and its interactions, shaping the world around us. Its profound impact is visible in every sphere of life, from pharmaceuticals and materials to energy and food production. Chemistry empowers us to decipher the mysteries of nature and harness its potential to create solutions for global challenges. It is a science that fuels progress, fosters innovation, and promises a better world for future generations."
Replace-Exact $sumOld $sumNew

# ---------------------------------------------------------------------
# 6. Add a new empty paragraph at the end of the document body
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 7. Font fix-up: TimesNewToman -> Times New Roman everywhere
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim().Length -gt 0) {
        $p.Range.Font.Name = "Times New Roman"
    }
}

Write-Output "done"
